$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price/volume snapshot (GitHub Actions scrape update).
# Cells use plain text (NumberFormat "@" forces Excel to keep numeric-looking
# strings like "7.50" / "494.40" / "1.00" as text instead of coercing them
# to numbers and stripping trailing zeros).
# Rows 45/46 also swap contents: FirstDigitalUSD and FLOKI traded ranking
# positions in this update.
$ws.Range("D2").Value = '70.880.59'
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = '3.811.98'
$ws.Range("E3").Value = '  -1.12%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '709.57'
$ws.Range("E5").Value = '  +1.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.58'
$ws.Range("E6").Value = '  -1.67%  '
$ws.Range("D7").Value = '3.811.19'
$ws.Range("E7").Value = '  -1.10%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.522'
$ws.Range("E9").Value = '  -0.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("E10").Value = '  -1.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.50'
$ws.Range("E11").Value = '  +3.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.457'
$ws.Range("E12").Value = '  -0.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("E13").Value = '  -1.70%  '
$ws.Range("E14").Value = '  -1.09%  '
$ws.Range("D15").Value = '4.455.93'
$ws.Range("E15").Value = '  -1.10%  '
$ws.Range("D16").Value = '3.827.49'
$ws.Range("E16").Value = '  -0.72%  '
$ws.Range("D17").Value = '70.899.44'
$ws.Range("E17").Value = '  -0.49%  '
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.14'
$ws.Range("E19").Value = '  -1.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.33'
$ws.Range("E20").Value = '  -2.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '494.40'
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.68'
$ws.Range("E22").Value = '  -4.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.727'
$ws.Range("E23").Value = '  +0.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.48'
$ws.Range("E24").Value = '  -0.73%  '
$ws.Range("E25").Value = '  -0.84%  '
$ws.Range("E26").Value = '  -1.64%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.44'
$ws.Range("E27").Value = '  -1.71%  '
$ws.Range("D28").Value = '3.963.39'
$ws.Range("E28").Value = '  -1.16%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  -4.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.09'
$ws.Range("E31").Value = '  -3.13%  '
$ws.Range("E32").Value = '  -1.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.35'
$ws.Range("E33").Value = '  -3.71%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.13'
$ws.Range("E34").Value = '  -1.94%  '
$ws.Range("E35").Value = '  -3.28%  '
$ws.Range("E36").Value = '  -1.64%  '
$ws.Range("D37").Value = '3.781.31'
$ws.Range("E37").Value = '  -0.68%  '
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("E39").Value = '  -2.41%  '
$ws.Range("E40").Value = '  +0.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.32'
$ws.Range("E41").Value = '  -3.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.95'
$ws.Range("E42").Value = '  -1.84%  '
$ws.Range("E43").Value = '  -3.91%  '
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("B45").Value = 'FLOKI'
$ws.Range("C45").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000324'
$ws.Range("E45").Value = '  +5.36%  '
$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '165.74'
$ws.Range("E47").Value = '  +1.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.86'
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '425.46'
$ws.Range("E49").Value = '  +1.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.62'
$ws.Range("E50").Value = '  -0.47%  '
$ws.Range("E51").Value = '  -2.65%  '
